# Refresh the cryptos price/volume snapshot (coinranking.com scrape).
# Most D (Price) / E (Volume 1h) cells get new readings; rows 25/26 and
# 46/47 also swap rank order (PancakeSwap<->Litecoin, VeChain<->ApeXProtocol).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.888.90'
$ws.Range("E2").Value = '  -1.11%  '

# Row 3
$ws.Range("D3").Value = '3.575.43'
$ws.Range("E3").Value = '  -2.01%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''575.72'
$ws.Range("E5").Value = '  -3.07%  '

# Row 6
$ws.Range("D6").Value = '''187.42'
$ws.Range("E6").Value = '  -3.43%  '

# Row 7
$ws.Range("D7").Value = '3.571.17'
$ws.Range("E7").Value = '  -2.00%  '

# Row 8
$ws.Range("E8").Value = '  -3.56%  '

# Row 10
$ws.Range("E10").Value = '  +2.30%  '

# Row 11
$ws.Range("E11").Value = '  -3.44%  '

# Row 12
$ws.Range("D12").Value = '''54.97'
$ws.Range("E12").Value = '  -5.44%  '

# Row 13
$ws.Range("D13").Value = '''0.0000303'
$ws.Range("E13").Value = '  +3.67%  '

# Row 14
$ws.Range("D14").Value = '''9.57'
$ws.Range("E14").Value = '  -3.47%  '

# Row 15
$ws.Range("D15").Value = '4.151.57'
$ws.Range("E15").Value = '  -1.88%  '

# Row 16
$ws.Range("D16").Value = '''19.69'
$ws.Range("E16").Value = '  -1.81%  '

# Row 17
$ws.Range("D17").Value = '3.574.98'
$ws.Range("E17").Value = '  -2.04%  '

# Row 18
$ws.Range("D18").Value = '69.881.09'
$ws.Range("E18").Value = '  -1.11%  '

# Row 19
$ws.Range("D19").Value = '''12.57'
$ws.Range("E19").Value = '  -1.36%  '

# Row 20
$ws.Range("E20").Value = '  -0.51%  '

# Row 21
$ws.Range("E21").Value = '  -3.10%  '

# Row 22
$ws.Range("D22").Value = '''489.19'
$ws.Range("E22").Value = '  -0.04%  '

# Row 23
$ws.Range("D23").Value = '''19.30'
$ws.Range("E23").Value = '  +1.90%  '

# Row 24
$ws.Range("D24").Value = '''4.90'
$ws.Range("E24").Value = '  -7.77%  '

# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '''4.39'
$ws.Range("E25").Value = '  -3.07%  '

# Row 26
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '''95.55'
$ws.Range("E26").Value = '  +4.96%  '

# Row 27
$ws.Range("D27").Value = '''11.30'
$ws.Range("E27").Value = '  -2.12%  '

# Row 28
$ws.Range("E28").Value = '  -6.47%  '

# Row 29
$ws.Range("D29").Value = '''9.32'
$ws.Range("E29").Value = '  -2.77%  '

# Row 30
$ws.Range("D30").Value = '''31.78'
$ws.Range("E30").Value = '  -3.00%  '

# Row 31
$ws.Range("D31").Value = '''7.64'
$ws.Range("E31").Value = '  -2.40%  '

# Row 32
$ws.Range("D32").Value = '''67.05'
$ws.Range("E32").Value = '  +2.10%  '

# Row 33
$ws.Range("D33").Value = '''12.06'
$ws.Range("E33").Value = '  -1.46%  '

# Row 34
$ws.Range("E34").Value = '  -5.01%  '

# Row 35
$ws.Range("D35").Value = '''569.03'
$ws.Range("E35").Value = '  -9.26%  '

# Row 36
$ws.Range("D36").Value = '''3.18'
$ws.Range("E36").Value = '  +12.49%  '

# Row 37
$ws.Range("D37").Value = '''38.72'
$ws.Range("E37").Value = '  -4.71%  '

# Row 38
$ws.Range("E38").Value = '  -0.01%  '

# Row 39
$ws.Range("E39").Value = '  -3.33%  '

# Row 40
$ws.Range("D40").Value = '''0.396'
$ws.Range("E40").Value = '  -3.47%  '

# Row 41
$ws.Range("E41").Value = '  +10.12%  '

# Row 42
$ws.Range("D42").Value = '''3.54'
$ws.Range("E42").Value = '  -0.81%  '

# Row 43
$ws.Range("E43").Value = '  -7.73%  '

# Row 44
$ws.Range("D44").Value = '3.266.61'
$ws.Range("E44").Value = '  -1.00%  '

# Row 45
$ws.Range("D45").Value = '''3.02'
$ws.Range("E45").Value = '  -4.14%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0438'
$ws.Range("E46").Value = '  -3.20%  '

# Row 47
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '''3.43'
$ws.Range("E47").Value = '  +3.95%  '

# Row 48
$ws.Range("D48").Value = '''9.68'
$ws.Range("E48").Value = '  +4.23%  '

# Row 49
$ws.Range("E49").Value = '  -1.84%  '

# Row 50
$ws.Range("D50").Value = '''0.999'
$ws.Range("E50").Value = '  -0.02%  '

# Row 51
$ws.Range("E51").Value = '  -3.73%  '
